# Generate Report for Handback
# Updates timestamps (and a priority flag) recorded by the handback status report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: refresh "Latest HO Xliff Generate Date" for the
#     8bb4e015-... entry (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-13 02:18:48"
$wsOverview.Range("G4").Value = "2016-08-13 02:18:48"

# --- zh-cn sheet: priority changed from human translation (ht) to
#     machine translation (mt), and handoff/handback timestamps refreshed
#     for the 8bb4e015-... entry (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-13 02:18:40"
$wsZhCn.Range("H4").Value = "2016-08-13 02:18:40"
$wsZhCn.Range("K3").Value = "2016-08-13 02:19:13"
$wsZhCn.Range("K4").Value = "2016-08-13 02:19:13"

# --- de-de sheet: handoff timestamp mirrors the Overview date, handback
#     timestamp refreshed for the 8bb4e015-... entry (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-13 02:18:48"
$wsDeDe.Range("H4").Value = "2016-08-13 02:18:48"
$wsDeDe.Range("K3").Value = "2016-08-13 02:19:22"
$wsDeDe.Range("K4").Value = "2016-08-13 02:19:22"
